$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block updates (B2, B5, B6, B7) ---
$ws.Range('B2').Value = "'2025-07-11"
$ws.Range('B2').Style = 'Normal'

$ws.Range('B5').Value = '
    • Studies published in English, peer-reviewed journals
    • About leptin and Alzheimer’s
    • Relevant papers available as full text
    • Randomized control trials 
    '
$ws.Rows(5).AutoFit()  # undo engine auto row-height bump from the embedded line breaks

$ws.Range('B6').Value = '40% of total quota selected for tranche'
$ws.Range('B7').Value = 'Randomized control trials'

# --- Row data updates (rows 12-20): Title/Year/Authors/Journal/Publisher/Summary ---
# Row 12
$ws.Range('B12').Value = 'Evolving cognition of the JAK-STAT signaling pathway: autoimmune disorders and cancer.'
$ws.Range('C12').Value = "'2023"
$ws.Range('C12').Style = 'Normal'
$ws.Range('D12').Value = 'Xue C, Yao Q, Gu X, Shi Q, Yuan X, Chu Q, Bao Z, Lu J, Li L'
$ws.Range('E12').Value = 'Signal transduction and targeted therapy'
$ws.Range('F12').Value = 'Signal Transduct Target Ther'
$ws.Range('G12').Value = 'The Janus kinase (JAK) signal transducer and activator of transcription (JAK-STAT) pathway is an evolutionarily conserved mechanism of transmembrane signal transduction that enables cells to communicate with the exterior environment. Various cytokines, interferons, growth factors, and other specific molecules activate JAK-STAT signaling to drive a series of physiological and pathological processes, including proliferation, metabolism, immune response, inflammation, and malignancy. Dysregulated JAK-STAT signaling and related genetic mutations are strongly associated with immune activation and cancer progression. Insights into the structures and functions of the JAK-STAT pathway have led to the development and approval of diverse drugs for the clinical treatment of diseases. Currently, drugs have been developed to mainly target the JAK-STAT pathway and are commonly divided into three subtypes: cytokine or receptor antibodies, JAK inhibitors, and STAT inhibitors. And novel agents also continue to be developed and tested in preclinical and clinical studies. The effectiveness and safety of each kind of drug also warrant further scientific trials before put into being clinical applications. Here, we review the current understanding of the fundamental composition and function of the JAK-STAT signaling pathway. We also discuss advancements in the understanding of JAK-STAT-related pathogenic mechanisms; targeted JAK-STAT therapies for various diseases, especially immune disorders, and cancers; newly developed JAK inhibitors; and current challenges and directions in the field.'

# Row 13
$ws.Range('B13').Value = 'The role of JAK/STAT signaling pathway and its inhibitors in diseases.'
$ws.Range('C13').Value = "'2020"
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').Value = 'Xin P, Xu X, Deng C, Liu S, Wang Y, Zhou X, Ma H, Wei D, Sun S'
$ws.Range('E13').Value = 'International immunopharmacology'
$ws.Range('F13').Value = 'Int Immunopharmacol'
$ws.Range('G13').Value = 'The JAK/STAT signaling pathway is an universally expressed intracellular signal transduction pathway and involved in many crucial biological processes, including cell proliferation, differentiation, apoptosis, and immune regulation. It provides a direct mechanism for extracellular factors-regulated gene expression. Current researches on this pathway have been focusing on the inflammatory and neoplastic diseases and related drug. The mechanism of JAK/STAT signaling is relatively simple. However, the biological consequences of the pathway are complicated due to its crosstalk with other signaling pathways. In addition, there is increasing evidence indicates that the persistent activation of JAK/STAT signaling pathway is closely related to many immune and inflammatory diseases, yet the specific mechanism remains unclear. Therefore, it is necessary to study the detailed mechanisms of JAK/STAT signaling in disease formation to provide critical reference for clinical treatments of the diseases. In this review, we focus on the structure of JAKs and STATs, the JAK/STAT signaling pathway and its negative regulators, the associated diseases, and the JAK inhibitors for the clinical therapy.'

# Row 14
$ws.Range('B14').Value = 'The JAK-STAT pathway at 30: Much learned, much more to do.'
$ws.Range('C14').Value = "'2022"
$ws.Range('C14').Style = 'Normal'
$ws.Range('D14').Value = 'Philips RL, Wang Y, Cheon H, Kanno Y, Gadina M, Sartorelli V, Horvath CM, Darnell JE, Stark GR, O''Shea JJ'
$ws.Range('E14').Value = 'Cell'
$ws.Range('F14').Value = 'Cell'
$ws.Range('G14').Value = 'The discovery of the Janus kinase (JAK)-signal transducer and activator of transcription (STAT) pathway arose from investigations of how cells respond to interferons (IFNs), revealing a paradigm in cell signaling conserved from slime molds to mammals. These discoveries revealed mechanisms underlying rapid gene expression mediated by a wide variety of extracellular polypeptides including cytokines, interleukins, and related factors. This knowledge has provided numerous insights into human disease, from immune deficiencies to cancer, and was rapidly translated to new drugs for autoimmune, allergic, and infectious diseases, including COVID-19. Despite these advances, major challenges and opportunities remain.'

# Row 15
$ws.Range('B15').Value = 'JAK-STAT signaling pathway in the pathogenesis of atopic dermatitis: An updated review.'
$ws.Range('C15').Value = "'2022"
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').Value = 'Huang IH, Chung WH, Wu PC, Chen CB'
$ws.Range('E15').Value = 'Frontiers in immunology'
$ws.Range('F15').Value = 'Front Immunol'
$ws.Range('G15').Value = 'Atopic dermatitis (AD) is a chronic, inflammatory, pruritic form of dermatosis with heterogeneous manifestations that can substantially affect patients'' quality of life. AD has a complex pathogenesis, making treatment challenging for dermatologists. The Janus kinase (JAK)-signal transducer and activator of transcription (STAT) pathway plays a central role in modulating multiple immune axes involved in the immunopathogenesis of AD. In particular, Th2 cytokines, including interleukin (IL)-4, IL-5, IL-13, IL-31, and thymic stromal lymphopoietin, which contribute to the symptoms of chronic inflammation and pruritus in AD, are mediated by JAK-STAT signal transduction. Furthermore, JAK-STAT is involved in the regulation of the epidermal barrier and the modulation of peripheral nerves related to the transduction of pruritus. Targeting the JAK-STAT pathway may attenuate these signals and show clinical efficacy through the suppression of various immune pathways associated with AD. Topical and oral JAK inhibitors with variable selectivity have emerged as promising therapeutic options for AD. Notably, topical ruxolitinib, oral upadacitinib, and oral abrocitinib were approved by the U.S. Food and Drug Administration for treating patients with AD. Accordingly, the present study reviewed the role of JAK-STAT pathways in the pathogenesis of AD and explored updated applications of JAK inhibitors in treating AD.'

# Row 16
$ws.Range('B16').Value = 'The molecular details of cytokine signaling via the JAK/STAT pathway.'
$ws.Range('C16').Value = "'2018"
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').Value = 'Morris R, Kershaw NJ, Babon JJ'
$ws.Range('E16').Value = 'Protein science : a publication of the Protein Society'
$ws.Range('F16').Value = 'Protein Sci'
$ws.Range('G16').Value = 'More than 50 cytokines signal via the JAK/STAT pathway to orchestrate hematopoiesis, induce inflammation and control the immune response. Cytokines are secreted glycoproteins that act as intercellular messengers, inducing proliferation, differentiation, growth, or apoptosis of their target cells. They act by binding to specific receptors on the surface of target cells and switching on a phosphotyrosine-based intracellular signaling cascade initiated by kinases then propagated and effected by SH2 domain-containing transcription factors. As cytokine signaling is proliferative and often inflammatory, it is tightly regulated in terms of both amplitude and duration. Here we review molecular details of the cytokine-induced signaling cascade and describe the architectures of the proteins involved, including the receptors, kinases, and transcription factors that initiate and propagate signaling and the regulatory proteins that control it.'

# Row 17
$ws.Range('B17').Value = 'Pathogenetic insights from the treatment of rheumatoid arthritis.'
$ws.Range('C17').Value = "'2017"
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').Value = 'McInnes IB, Schett G'
$ws.Range('E17').Value = 'Lancet (London, England)'
$ws.Range('F17').Value = 'Lancet'
$ws.Range('G17').Value = 'Rheumatoid arthritis is a chronic autoimmune disease that causes progressive articular damage, functional loss, and comorbidity. The development of effective biologics and small-molecule kinase inhibitors in the past two decades has substantially improved clinical outcomes. Just as understanding of pathogenesis has led in large part to the development of drugs, so have mode-of-action studies of these specific immune-targeted agents revealed which immune pathways drive articular inflammation and related comorbidities. Cytokine inhibitors have definitively proven a critical role for tumour necrosis factor α and interleukin 6 in disease pathogenesis and possibly also for granulocyte-macrophage colony-stimulating factor. More recently, clinical trials with Janus kinase (JAK) inhibitors have shown that cytokine receptors that signal through the JAK/STAT signalling pathway are important for disease, informing the pathogenetic function of additional cytokines (such as the interferons). Finally, successful use of costimulatory blockade and B-cell depletion in the clinic has revealed that the adaptive immune response and the downstream events initiated by these cells participate directly in synovial inflammation. Taken together, it becomes apparent that understanding the effects of specific immune interventions can elucidate definitive molecular or cellular nodes that are essential to maintain complex inflammatory networks that subserve diseases like rheumatoid arthritis.'

# Row 18
$ws.Range('B18').Value = 'Principles of interleukin (IL)-6-type cytokine signalling and its regulation.'
$ws.Range('C18').Value = "'2003"
$ws.Range('C18').Style = 'Normal'
$ws.Range('D18').Value = 'Heinrich PC, Behrmann I, Haan S, Hermanns HM, Müller-Newen G, Schaper F'
$ws.Range('E18').Value = 'The Biochemical journal'
$ws.Range('F18').Value = 'Biochem J'
$ws.Range('G18').Value = 'The IL (interleukin)-6-type cytokines IL-6, IL-11, LIF (leukaemia inhibitory factor), OSM (oncostatin M), ciliary neurotrophic factor, cardiotrophin-1 and cardiotrophin-like cytokine are an important family of mediators involved in the regulation of the acute-phase response to injury and infection. Besides their functions in inflammation and the immune response, these cytokines play also a crucial role in haematopoiesis, liver and neuronal regeneration, embryonal development and fertility. Dysregulation of IL-6-type cytokine signalling contributes to the onset and maintenance of several diseases, such as rheumatoid arthritis, inflammatory bowel disease, osteoporosis, multiple sclerosis and various types of cancer (e.g. multiple myeloma and prostate cancer). IL-6-type cytokines exert their action via the signal transducers gp (glycoprotein) 130, LIF receptor and OSM receptor leading to the activation of the JAK/STAT (Janus kinase/signal transducer and activator of transcription) and MAPK (mitogen-activated protein kinase) cascades. This review focuses on recent progress in the understanding of the molecular mechanisms of IL-6-type cytokine signal transduction. Emphasis is put on the termination and modulation of the JAK/STAT signalling pathway mediated by tyrosine phosphatases, the SOCS (suppressor of cytokine signalling) feedback inhibitors and PIAS (protein inhibitor of activated STAT) proteins. Also the cross-talk between the JAK/STAT pathway with other signalling cascades is discussed.'

# Row 19
$ws.Range('B19').Value = 'JAK-STAT pathway targeting for the treatment of inflammatory bowel disease.'
$ws.Range('C19').Value = "'2020"
$ws.Range('C19').Style = 'Normal'
$ws.Range('D19').Value = 'Salas A, Hernandez-Rocha C, Duijvestein M, Faubion W, McGovern D, Vermeire S, Vetrano S, Vande Casteele N'
$ws.Range('E19').Value = 'Nature reviews. Gastroenterology & hepatology'
$ws.Range('F19').Value = 'Nat Rev Gastroenterol Hepatol'
$ws.Range('G19').Value = 'Cytokines are involved in intestinal homeostasis and pathological processes associated with inflammatory bowel disease (IBD). The biological effects of cytokines, including several involved in the pathology of Crohn''s disease and ulcerative colitis, occur as a result of receptor-mediated signalling through the Janus kinase (JAK) and signal transducer and activator of transcription (STAT) DNA-binding families of proteins. Although therapies targeting cytokines have revolutionized IBD therapy, they have historically targeted individual cytokines, and an unmet medical need exists for patients who do not respond to or lose response to these treatments. Several small-molecule inhibitors of JAKs that have the potential to affect multiple pro-inflammatory cytokine-dependent pathways are in clinical development for the treatment of IBD, with one agent, tofacitinib, already approved for ulcerative colitis and several other agents with demonstrated efficacy in early phase trials. This Review describes the current understanding of JAK-STAT signalling in intestinal homeostasis and disease and the rationale for targeting this pathway as a treatment for IBD. The available evidence for the efficacy, safety and pharmacokinetics of JAK inhibitors in IBD as well as the potential approaches to optimize treatment with these agents, such as localized delivery or combination therapy, are also discussed.'

# Row 20
$ws.Range('B20').Value = 'JAK inhibitor: Introduction.'
$ws.Range('C20').Value = "'2023"
$ws.Range('C20').Style = 'Normal'
$ws.Range('D20').Value = 'Raychaudhuri SP, Raychaudhuri SK'
$ws.Range('E20').Value = 'Indian journal of dermatology, venereology and leprology'
$ws.Range('F20').Value = 'Indian J Dermatol Venereol Leprol'
$ws.Range('G20').Value = 'The Janus kinase (JAK)-signal transducer and activator of transcription (STAT) pathway is a key regulatory signaling system for cellular proliferation, differentiation, and apoptosis. In addition, JAK signaling pathway plays critical roles in orchestrating immune response through its interactions with the cytokine receptors and the transcriptions factors. Several key cytokines use JAK-STAT signaling proteins to transduce intra-cellular signals which are involved in the pathogenesis of autoimmune and inflammatory diseases such as in psoriatic disease (psoriasis, psoriatic arthritis), atopic dermatitis, alopecia areata, vitiligo, rheumatoid arthritis, ankylosing spondylitis, lupus erythematosus, Sjogren''s syndrome, and other autoimmune diseases. In recent years, understandings of the molecular mechanisms of JAK-STAT pathway in the inflammatory proliferative cascades of autoimmune diseases has led to the development of JAK inhibitors and has opened a new dimension for the treatment of systemic and cutaneous inflammatory diseases. In this symposium we have provided a broad perspective on the use of Janus kinase inhibitors in cutaneous autoimmune diseases.'

# --- New row 21 ---
$ws.Range('A21').Value = 10
$ws.Range('B21').Value = 'The JAK-STAT pathway: impact on human disease and therapeutic intervention.'
$ws.Range('C21').Value = "'2015"
$ws.Range('C21').Style = 'Normal'
$ws.Range('D21').Value = 'O''Shea JJ, Schwartz DM, Villarino AV, Gadina M, McInnes IB, Laurence A'
$ws.Range('E21').Value = 'Annual review of medicine'
$ws.Range('F21').Value = 'Annu Rev Med'
$ws.Range('G21').Value = 'The Janus kinase (JAK)-signal transducer of activators of transcription (STAT) pathway is now recognized as an evolutionarily conserved signaling pathway employed by diverse cytokines, interferons, growth factors, and related molecules. This pathway provides an elegant and remarkably straightforward mechanism whereby extracellular factors control gene expression. It thus serves as a fundamental paradigm for how cells sense environmental cues and interpret these signals to regulate cell growth and differentiation. Genetic mutations and polymorphisms are functionally relevant to a variety of human diseases, especially cancer and immune-related conditions. The clinical relevance of the pathway has been confirmed by the emergence of a new class of therapeutics that targets JAKs.'

# --- Extend AutoFilter to cover the new row, and update the hidden _FilterDatabase name ---
$ws.AutoFilterMode = $false
$ws.Range('A11:I21').AutoFilter()
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Search Results'!`$A`$11:`$I`$21"
    }
}
